$d = $word.ActiveDocument

# Paragraph 8 (the final, empty paragraph) currently holds the "_GoBack"
# bookmark. Remove it entirely -- it (and its bookmark) should disappear.
$last = $d.Paragraphs.Item(8)
$last.Range.Delete()

# Paragraph 5 held the "Display function does not run every time it is
# called?" sentence. Clear that text and drop the "_GoBack" bookmark onto
# the now-empty paragraph instead.
$p5 = $d.Paragraphs.Item(5)
$r = $p5.Range
[void]$r.MoveEnd(1, -1)
$r.Text = ""
[void]$d.Bookmarks.Add("_GoBack", $r)
